# Add "Norway" and "Poland" market test-data sheets, cloned from the
# existing "Croatia" sheet (same layout/column widths/no custom row
# heights as the template used for these two new country sheets), then
# fill in the market name / NGC code for each and leave "Norway" as the
# active tab.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Croatia")

# --- Norway ------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
# Set the NGC code (B4) before the market name (B2) so the shared-string
# table picks up the codes first, matching how the existing country
# sheets were authored.
$norway.Range("B4").Value = "NGC-2931/T3071/T3070/T3072"
$norway.Range("B2").Value = "Norway Market"

# --- Poland --------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet2)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/3036/T3037/T3038"
$poland.Range("B2").Value = "Poland Market"

# Leave "Norway" as the selected/active sheet.
$norway.Activate()
